$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append two new bullet items to the "Vrijdag 15 november 2019:"
#    list (numId 4), right after the existing "Wireframes gemaakt"
#    item: "Technisch ontwerp gemaakt" and "Functioneel ontwerp
#    begonnen". Inserting a paragraph right after an existing list
#    item inherits that item's paragraph/run formatting (ListParagraph
#    style, numId 4, Arial 12pt), so no extra formatting needs to be
#    applied explicitly.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)                      # wdCollapseEnd
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.InsertBefore("Technisch ontwerp gemaakt")

$newPara1.Range.Collapse(0)         # wdCollapseEnd
$newPara1.Range.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.InsertBefore("Functioneel ontwerp begonnen")

# ------------------------------------------------------------------
# 2. The document carries a (hidden) "_GoBack" bookmark that marks the
#    last edit location. Originally it sat right after "Wireframes
#    gemaakt" (now no longer the last edited text); move it so it
#    sits right after the newly typed "Functioneel ontwerp begonnen",
#    matching a normal editing session where the last thing typed was
#    that final bullet.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Add a one-character placeholder at the very end of the paragraph so
# that the position just before it is not itself adjacent to the
# paragraph mark, then collapse-add the bookmark there and remove the
# placeholder again. (Adding a zero-length bookmark whose position is
# immediately at a paragraph mark is unreliable, so this keeps the
# bookmark insertion point safely inside the paragraph's text while
# it is created.)
$endR = $finalPara.Range
$endR.Collapse(0)                   # wdCollapseEnd
$endR.InsertBefore("X")

$markerPos = $finalPara.Range.End - 2
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($markerPos, $markerPos + 1)
$placeholderRange.Delete()
